$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A6").Value = 148
$ws.Range("A7").Value = 149
$ws.Range("A8").Value = 150
